# Apply the "update class diagram + usecase + physical diagram 20160411" edit:
# fill in rows 61 and 62 of Sheet1 (previously only had the "No" column filled
# in) with Function / Description / Status data, and move the sheet's
# viewport/selection down to show the newly-populated rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61: Server Index / Lỗi filter bảng list servers / Pending
$ws.Range("B61").Value = "Server Index"
$ws.Range("C61").Value = "Lỗi filter bảng list servers"
$ws.Range("E61").Value = "Pending"

# Row 62: Ca kiếp / liệt kê hết staff ... / Pending
$ws.Range("B62").Value = "Ca kiếp"
$ws.Range("C62").Value = "liệt kê hết staff trong hệ thống ra --> sửa lại trường hợp nó cho đổi ca trực"
$ws.Range("E62").Value = "Pending"

# Row 62's description wraps onto a second line, so the row grows taller
# (matches the other two-line wrapped description rows in this sheet).
$ws.Rows.Item(62).RowHeight = 30

# Move the visible window down and reselect the new active cell, same as the
# author's view state after editing these rows.
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E63").Select()
